$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 133: "Big Brush, Big Dreams" / "Ginseng Angle Brush"
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -60120

# ALC row 138: "All-night Crafting" / "Cunning Craftsman's Tisane"
$ws.Range("H138").Value = 7237634
$ws.Range("I138").Value = 6498414
$ws.Range("J138").Value = 7357214
$ws.Range("K138").Value = 19495242
$ws.Range("L138").Value = 22071642
$ws.Range("M138").Value = -19490102
$ws.Range("N138").Value = -22081922

$ws = $wb.Worksheets.Item("ARM")
# ARM row 4: "Eyes Bigger than the Plate" / "Bronze Plate"
$ws.Range("H4").Value = 1487.375
$ws.Range("I4").Value = 1487.375
$ws.Range("K4").Value = 1487.375
$ws.Range("M4").Value = -1371.375

# ARM row 44: "Very Slow Array" / "Mythril Plate"
$ws.Range("H44").Value = 30000
$ws.Range("J44").Value = 30000
$ws.Range("L44").Value = 30000
$ws.Range("N44").Value = -30976

# ARM row 55: "Employee Retention" / "Mythril Elmo"
$ws.Range("H55").Value = 7355.5557
$ws.Range("I55").Value = 7000
$ws.Range("K55").Value = 7000
$ws.Range("M55").Value = -6685

# ARM row 61: "Dealing with the Tough Stuff" / "Cobalt Ingot"
$ws.Range("H61").Value = 3044.95
$ws.Range("I61").Value = 2648.9333
$ws.Range("J61").Value = 4233
$ws.Range("K61").Value = 2648.9333
$ws.Range("L61").Value = 4233
$ws.Range("M61").Value = -2436.9333
$ws.Range("N61").Value = -4657

# ARM row 132: "Don't Bore Me, Ore Me" / "Mountain Chromite Ingot"
$ws.Range("H132").Value = 1972.3448
$ws.Range("I132").Value = 1872.8085
$ws.Range("J132").Value = 2397.6365
$ws.Range("K132").Value = 5618.4255
$ws.Range("L132").Value = 7192.9095
$ws.Range("M132").Value = -3088.4255
$ws.Range("N132").Value = -12252.9095

# ARM row 136: "Metal with Mettle" / "Cobalt Tungsten Ingot"
$ws.Range("H136").Value = 3044.95
$ws.Range("I136").Value = 2648.9333
$ws.Range("J136").Value = 4233
$ws.Range("K136").Value = 7946.7999
$ws.Range("L136").Value = 12699
$ws.Range("M136").Value = -5396.7999
$ws.Range("N136").Value = -17799

# ARM row 138: "Don't Ask about the Rivets" / "Titanium Gold Helm of Casting"
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").ClearContents()
$ws.Range("N138").Value = 0

# ARM row 139: "Backing up My Words" / "Titanium Gold Thornplate of Fending"
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("N139").Value = 0

$ws = $wb.Worksheets.Item("BSM")
# BSM row 105: "Ingot to Wing It" / "Molybdenum Ingot"
$ws.Range("H105").Value = 2979.6829
$ws.Range("I105").Value = 2782.2334
$ws.Range("J105").Value = 3518.182
$ws.Range("K105").Value = 2782.2334
$ws.Range("L105").Value = 3518.182
$ws.Range("M105").Value = -1035.2334
$ws.Range("N105").Value = -7012.182

# BSM row 135: "Axes to the Maxes" / "Ruthenium War Axe"
$ws.Range("H135").Value = 54640
$ws.Range("J135").Value = 54640
$ws.Range("L135").Value = 54640
$ws.Range("N135").Value = -64780

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31: "Wall Not Found" / "Walnut Lumber"
$ws.Range("H31").Value = 5789.204
$ws.Range("I31").Value = 2988.2144
$ws.Range("J31").Value = 6909.6
$ws.Range("K31").Value = 2988.2144
$ws.Range("L31").Value = 6909.6
$ws.Range("M31").Value = -2693.2144
$ws.Range("N31").Value = -7499.6

# CRP row 34: "Armoires of the Rich and Famous" / "Walnut Lumber"
$ws.Range("H34").Value = 5789.204
$ws.Range("I34").Value = 2988.2144
$ws.Range("J34").Value = 6909.6
$ws.Range("K34").Value = 2988.2144
$ws.Range("L34").Value = 6909.6
$ws.Range("M34").Value = -2786.2144
$ws.Range("N34").Value = -7313.6

# CRP row 99: "O Pine" / "Pine Lumber"
$ws.Range("H99").Value = 19612284
$ws.Range("J99").Value = 55559640
$ws.Range("L99").Value = 55559640
$ws.Range("N99").Value = -55562636

# CRP row 126: "A Better Conductor" / "Red Pine Lumber"
$ws.Range("H126").Value = 19612284
$ws.Range("J126").Value = 55559640
$ws.Range("L126").Value = 166678920
$ws.Range("N126").Value = -166683860

$ws = $wb.Worksheets.Item("CUL")
# CUL row 5: "What a Sap" / "Maple Syrup"
$ws.Range("H5").Value = 2123.5715
$ws.Range("I5").Value = 891.06665
$ws.Range("J5").Value = 3545.6924
$ws.Range("K5").Value = 2673.19995
$ws.Range("L5").Value = 10637.0772
$ws.Range("M5").Value = -2561.19995
$ws.Range("N5").Value = -10861.0772

# CUL row 22: "A Total Nut Job" / "Walnut Bread"
$ws.Range("H22").Value = 1450
$ws.Range("J22").Value = 2500
$ws.Range("L22").Value = 7500
$ws.Range("N22").Value = -7838

# CUL row 27: "Brain Food" / "Walnut Bread"
$ws.Range("H27").Value = 1450
$ws.Range("J27").Value = 2500
$ws.Range("L27").Value = 7500
$ws.Range("N27").Value = -7704

# CUL row 113: "Can't Eat Just One" / "Night Vinegar"
$ws.Range("H113").Value = 6994146.5
$ws.Range("I113").Value = 800
$ws.Range("J113").Value = 7576925
$ws.Range("K113").Value = 2400
$ws.Range("L113").Value = 22730775
$ws.Range("M113").Value = -230
$ws.Range("N113").Value = -22735115

# CUL row 121: "A Cookie for Your Troubles" / "Coffee Biscuit"
$ws.Range("H121").Value = 1067.9445
$ws.Range("I121").Value = 228.33333
$ws.Range("J121").Value = 1487.75
$ws.Range("K121").Value = 684.99999
$ws.Range("L121").Value = 4463.25
$ws.Range("M121").Value = 625.00001
$ws.Range("N121").Value = -7083.25

# CUL row 135: "Not-so-secret Ingredient" / "Royal Maple Syrup"
$ws.Range("H135").Value = 2123.5715
$ws.Range("I135").Value = 891.06665
$ws.Range("J135").Value = 3545.6924
$ws.Range("K135").Value = 8019.59985
$ws.Range("L135").Value = 31911.2316
$ws.Range("M135").Value = -5484.59985
$ws.Range("N135").Value = -36981.2316

$ws = $wb.Worksheets.Item("GSM")
# GSM row 17: "Point of Honor" / "Amateur's Needle"
$ws.Range("H17").Value = 20000
$ws.Range("J17").Value = 5000
$ws.Range("L17").Value = 5000
$ws.Range("N17").Value = -5336

# GSM row 32: "Love in the Time of Umbra" / "Silver Ring"
$ws.Range("H32").Value = 29290
$ws.Range("J32").Value = 29290
$ws.Range("L32").Value = 29290
$ws.Range("N32").Value = -29882

# GSM row 42: "It's Only Love" / "Silver Choker"
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").ClearContents()
$ws.Range("N42").Value = 0

# GSM row 97: "If I'd a Koppranickel for Every Time..." / "Koppranickel Ingot"
$ws.Range("H97").Value = 1534.9231
$ws.Range("I97").Value = 1435.3889
$ws.Range("J97").Value = 1758.875
$ws.Range("K97").Value = 1435.3889
$ws.Range("L97").Value = 1758.875
$ws.Range("M97").Value = -939.3888999999999
$ws.Range("N97").Value = -2750.875

# GSM row 102: "Put the Metal to the Peddle" / "Durium Ingot"
$ws.Range("H102").Value = 8198.777
$ws.Range("I102").Value = 4331.5
$ws.Range("K102").Value = 4331.5
$ws.Range("M102").Value = -2709.5

# GSM row 115: "Unsung Generosity" / "Manasilver Choker"
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").ClearContents()
$ws.Range("N115").Value = 0

# GSM row 122: "Awarding Academic Excellence" / "Ametrine"
$ws.Range("H122").Value = 3315.1904
$ws.Range("I122").Value = 3279
$ws.Range("J122").Value = 3532.3333
$ws.Range("K122").Value = 9837
$ws.Range("L122").Value = 10596.9999
$ws.Range("M122").Value = -7387
$ws.Range("N122").Value = -15496.9999

$ws = $wb.Worksheets.Item("LTW")
# LTW row 132: "Tenets of Tanning" / "Silver Lobo Leather"
$ws.Range("H132").Value = 2983.279
$ws.Range("I132").Value = 2212.718
$ws.Range("K132").Value = 6638.154
$ws.Range("M132").Value = -4108.154

# LTW row 136: "Respect for Br'aax" / "Br'aax Leather"
$ws.Range("H136").Value = 7717.6875
$ws.Range("I136").Value = 2957.8
$ws.Range("J136").Value = 15650.833
$ws.Range("K136").Value = 8873.400000000001
$ws.Range("L136").Value = 46952.499
$ws.Range("M136").Value = -6323.400000000001
$ws.Range("N136").Value = -52052.499

$ws = $wb.Worksheets.Item("WVR")
# WVR row 4: "Not Cool Enough" / "Hempen Undershirt"
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").ClearContents()
$ws.Range("N4").Value = 0
